$wb = $excel.ActiveWorkbook

# The "Portugal" sheet is the template for the new "Slovakia" sheet: duplicate it,
# placing the copy immediately after Portugal (this becomes the newest worksheet part).
$portugal = $wb.Worksheets.Item("Portugal")
$portugal.Copy([System.Reflection.Missing]::Value, $portugal)
$slovakia = $wb.Worksheets.Item($portugal.Index + 1)
$slovakia.Name = "Slovakia"

# Fill in the Slovakia-specific market name and Jira ticket reference.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3236/T3235/T3234/T3237/T3228"

# The rows holding those two cells re-wrap/auto-size now that the text differs.
$slovakia.Rows.Item(3).AutoFit() | Out-Null
$slovakia.Rows.Item(4).AutoFit() | Out-Null

# Reset the view back on Portugal (no longer the active tab) to a "whole sheet" selection.
$portugal.Cells.Select()

# Make the new Slovakia sheet the active tab/view.
$slovakia.Activate()
$slovakia.Range("B4").Select()
